$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Standard Error (SE) values for the "Ascidian_transect_percent_cover"
# dataset, column F, rows 2-130 (one data row per TP/Site/Transect_ID).
$newSE = @(0, 0, 2.33428552000155, 0.213437474581095, 0.690410505906933, 0.4, 0.213437474581095, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0.442216638714053, 0, 0, 0, 0.1, 0, 0, 0.1, 0, 0.213437474581095, 0, 0, 0.1, 0, 0, 0, 0, 0.1, 0, 0, 0.3, 0, 0, 0, 0, 0, 1, 0, 0.6, 0, 0, 0.1, 0, 0, 0, 0, 0.1, 1.18181818181818, 0.90921211313239, 0, 0, 0.1, 0.272727272727273, 0, 0.263479577203444, 0.2, 0, 0, 0, 0.1, 0.703167436990966, 0, 0, 0.1, 0.2, 0.221108319357027, 0.517472489875334, 0, 0, 0.517472489875334, 0, 1.03494497975067, 0, 0, 0, 0.2, 0, 0, 0, 0, 2.9152853872115, 1.43139403690559, 3.91918358845308, 0, 0, 0, 0, 0.1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0.152752523165195, 0, 0, 0, 0, 0, 3.13758576686669, 0, 0, 0.5, 0, 0.3, 0.2, 0, 0.2, 0, 0, 0, 0, 0, 0, 0, 0.1, 0, 0, 0)

$startRow = 2
for ($i = 0; $i -lt $newSE.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 6).Value = $newSE[$i]
}
